$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 38.1
$ws.Range("C4").Value = 23.68
$ws.Range("C5").Value = 4.51
$ws.Range("C6").Value = 1417
$ws.Range("C7").Value = 21495.6065
$ws.Range("C8").Value = 1050
$ws.Range("C9").Value = 474
$ws.Range("C10").Value = 8.6899
$ws.Range("C11").Value = 7.7099
$ws.Range("C12").Value = 110654.9954
$ws.Range("C13").Value = 37410.1999
$ws.Range("C14").Value = 642.1699
$ws.Range("C15").Value = 514.8644
$ws.Range("C16").Value = 324.6861
$ws.Range("C17").Value = 3.8203
$ws.Range("C18").Value = 13.1989
$ws.Range("C19").Value = 0.0455

$ws.Range("C8").Select()
